$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 'QUIQUIA MALLQUI, CYNTHIA ANGELLINE'
$ws.Range("B2").Value = '06:30'
$ws.Range("C2").Value = '10:15'
$ws.Range("D2").Value = 'QUIQUIA MALLQUI, CYNTHIA ANGELLINE'
$ws.Range("E2").Value = '10:15'
$ws.Range("F2").Value = '06:30'

$ws.Range("A3").Value = 'VEGA RIVAS, ANDREA FERNANDA'
$ws.Range("B3").Value = '07:00'
$ws.Range("C3").Value = '10:45'
$ws.Range("D3").Value = 'VEGA RIVAS, ANDREA FERNANDA'
$ws.Range("E3").Value = '10:45'
$ws.Range("F3").Value = '07:00'

$ws.Range("A4").Value = 'YOVERA ROBLES, VICTOR EDUARDO'
$ws.Range("B4").Value = '08:00'
$ws.Range("C4").Value = '11:45'
$ws.Range("D4").Value = 'YOVERA ROBLES, VICTOR EDUARDO'
$ws.Range("E4").Value = '11:45'
$ws.Range("F4").Value = '08:00'

$ws.Range("A5").Value = 'MEZA PEREZ, JUAN CRISTOFER'
$ws.Range("B5").Value = '08:00'
$ws.Range("C5").Value = '11:45'
$ws.Range("D5").Value = 'MEZA PEREZ, JUAN CRISTOFER'
$ws.Range("E5").Value = '11:45'
$ws.Range("F5").Value = '08:00'

$ws.Range("A6").Value = 'MARTINEZ PAZ, ROCIO ESPERANZA'
$ws.Range("B6").Value = '09:00'
$ws.Range("C6").Value = '18:00'
$ws.Range("D6").Value = 'CARHUARICRA ESPINOZA, FIORELLA NICOLL'
$ws.Range("E6").Value = '12:45'
$ws.Range("F6").Value = '09:00'

$ws.Range("A7").Value = 'CARHUARICRA ESPINOZA, FIORELLA NICOLL'
$ws.Range("B7").Value = '09:00'
$ws.Range("C7").Value = '12:45'
$ws.Range("D7").Value = 'YANQUI BRAVO, MIRIAN LUZ'
$ws.Range("E7").Value = '13:00'
$ws.Range("F7").Value = '09:15'

$ws.Range("A8").Value = 'YANQUI BRAVO, MIRIAN LUZ'
$ws.Range("B8").Value = '09:15'
$ws.Range("C8").Value = '13:00'
$ws.Range("D8").Value = 'CUSI QUISPE, ANDREA ESTEFANY'
$ws.Range("E8").Value = '13:30'
$ws.Range("F8").Value = '09:45'

$ws.Range("A9").Value = 'CUSI QUISPE, ANDREA ESTEFANY'
$ws.Range("B9").Value = '09:45'
$ws.Range("C9").Value = '13:30'
$ws.Range("D9").Value = 'GOMEZ ALBINO, IDALIA GIMENA'
$ws.Range("E9").Value = '13:45'
$ws.Range("F9").Value = '10:00'

$ws.Range("A10").Value = 'GOMEZ ALBINO, IDALIA GIMENA'
$ws.Range("B10").Value = '10:00'
$ws.Range("C10").Value = '13:45'
$ws.Range("D10").Value = 'MONTEZUMA DEJO, EVELYN BRUNELLA'
$ws.Range("E10").Value = '14:00'
$ws.Range("F10").Value = '10:15'

$ws.Range("A11").Value = 'MONTEZUMA DEJO, EVELYN BRUNELLA'
$ws.Range("B11").Value = '10:15'
$ws.Range("C11").Value = '14:00'
$ws.Range("D11").Value = 'HUAYANAY VELASCO, ATHINA'
$ws.Range("E11").Value = '14:00'
$ws.Range("F11").Value = '10:15'

$ws.Range("A12").Value = 'HUAYANAY VELASCO, ATHINA'
$ws.Range("B12").Value = '10:15'
$ws.Range("C12").Value = '14:00'
$ws.Range("D12").Value = 'MORENO CANCHANYA, ROSMERY'
$ws.Range("E12").Value = '14:15'
$ws.Range("F12").Value = '10:30'

$ws.Range("A13").Value = 'MORENO CANCHANYA, ROSMERY'
$ws.Range("B13").Value = '10:30'
$ws.Range("C13").Value = '14:15'
$ws.Range("D13").Value = 'RUIZ SANTOS, CIELO CRISTHINA'
$ws.Range("E13").Value = '14:30'
$ws.Range("F13").Value = '10:45'

$ws.Range("A14").Value = 'RUIZ SANTOS, CIELO CRISTHINA'
$ws.Range("B14").Value = '10:45'
$ws.Range("C14").Value = '14:30'
$ws.Range("D14").Value = 'VEGA CARDENAS, ANGELICA LOURDES'
$ws.Range("E14").Value = '14:45'
$ws.Range("F14").Value = '11:00'

$ws.Range("A15").Value = 'VEGA CARDENAS, ANGELICA LOURDES'
$ws.Range("B15").Value = '11:00'
$ws.Range("C15").Value = '14:45'
$ws.Range("D15").Value = 'HUAMAN HUAMANI, ALEXIS JAVIER'
$ws.Range("E15").Value = '14:45'
$ws.Range("F15").Value = '11:00'

$ws.Range("A16").Value = 'HUAMAN HUAMANI, ALEXIS JAVIER'
$ws.Range("B16").Value = '11:00'
$ws.Range("C16").Value = '14:45'
$ws.Range("D16").Value = 'Del Aguila Murayari, Darla'
$ws.Range("E16").Value = '14:45'
$ws.Range("F16").Value = '11:00'

$ws.Range("A17").Value = 'Del Aguila Murayari, Darla'
$ws.Range("B17").Value = '11:00'
$ws.Range("C17").Value = '14:45'
$ws.Range("D17").Value = 'QUISPE MONDRAGÓN, JUAN ALFONSO'
$ws.Range("E17").Value = '15:00'
$ws.Range("F17").Value = '11:15'

$ws.Range("A18").Value = 'QUISPE MONDRAGÓN, JUAN ALFONSO'
$ws.Range("B18").Value = '11:15'
$ws.Range("C18").Value = '15:00'
$ws.Range("D18").Value = 'RIVERA CARREÑO, DIANA DESIRÉE'
$ws.Range("E18").Value = '17:45'
$ws.Range("F18").Value = '14:00'

$ws.Range("A19").Value = 'SOTELO GONZALES, CAMILA SOFÍA'
$ws.Range("B19").Value = '12:30'
$ws.Range("C19").Value = '21:30'
$ws.Range("D19").Value = 'MARTINEZ PAZ, ROCIO ESPERANZA'
$ws.Range("E19").Value = '18:00'
$ws.Range("F19").Value = '09:00'

$ws.Range("A20").Value = 'RIVERA CARREÑO, DIANA DESIRÉE'
$ws.Range("B20").Value = '14:00'
$ws.Range("C20").Value = '17:45'
$ws.Range("D20").Value = 'SUAREZ JARA, YENNIFER YUSSARA'
$ws.Range("E20").Value = '18:15'
$ws.Range("F20").Value = '14:30'

$ws.Range("A21").Value = 'SUAREZ JARA, YENNIFER YUSSARA'
$ws.Range("B21").Value = '14:30'
$ws.Range("C21").Value = '18:15'
$ws.Range("D21").Value = 'VILCAPOMA CHILIN, JULISSA JAZMIN'
$ws.Range("E21").Value = '18:15'
$ws.Range("F21").Value = '14:30'

$ws.Range("A22").Value = 'VILCAPOMA CHILIN, JULISSA JAZMIN'
$ws.Range("B22").Value = '14:30'
$ws.Range("C22").Value = '18:15'
$ws.Range("D22").Value = 'MARTICORENA LOPEZ, DAVID CARLOS'
$ws.Range("E22").Value = '18:30'
$ws.Range("F22").Value = '14:45'

$ws.Range("A23").Value = 'MARTICORENA LOPEZ, DAVID CARLOS'
$ws.Range("B23").Value = '14:45'
$ws.Range("C23").Value = '18:30'
$ws.Range("D23").Value = 'MEDINA MARCELO, NAOMI ARIADNA'
$ws.Range("E23").Value = '18:45'
$ws.Range("F23").Value = '15:00'

$ws.Range("A24").Value = 'MEDINA MARCELO, NAOMI ARIADNA'
$ws.Range("B24").Value = '15:00'
$ws.Range("C24").Value = '18:45'
$ws.Range("D24").Value = 'MUÑOZ SOTOMAYOR, MIRIAN RAQUEL'
$ws.Range("E24").Value = '19:00'
$ws.Range("F24").Value = '15:15'

$ws.Range("A25").Value = 'MUÑOZ SOTOMAYOR, MIRIAN RAQUEL'
$ws.Range("B25").Value = '15:15'
$ws.Range("C25").Value = '19:00'
$ws.Range("D25").Value = 'VARGAS CASTRO, LOANA VICTORIA'
$ws.Range("E25").Value = '20:00'
$ws.Range("F25").Value = '16:15'

$ws.Range("A26").Value = 'VARGAS CASTRO, LOANA VICTORIA'
$ws.Range("B26").Value = '16:15'
$ws.Range("C26").Value = '20:00'
$ws.Range("D26").Value = 'TORRES RAZURI, JESUS GUSTAVO SANTIAGO'
$ws.Range("E26").Value = '20:15'
$ws.Range("F26").Value = '16:30'

$ws.Range("A27").Value = 'TORRES RAZURI, JESUS GUSTAVO SANTIAGO'
$ws.Range("B27").Value = '16:30'
$ws.Range("C27").Value = '20:15'
$ws.Range("D27").Value = 'YACILA GRANDEZ, RODRIGO ANDRE'
$ws.Range("E27").Value = '20:15'
$ws.Range("F27").Value = '16:30'

$ws.Range("A28").Value = 'YACILA GRANDEZ, RODRIGO ANDRE'
$ws.Range("B28").Value = '16:30'
$ws.Range("C28").Value = '20:15'
$ws.Range("D28").Value = 'CHIARA LIMA, AUGUSTO SEBASTIAN'
$ws.Range("E28").Value = '21:00'
$ws.Range("F28").Value = '17:15'

$ws.Range("A29").Value = 'CHIARA LIMA, AUGUSTO SEBASTIAN'
$ws.Range("B29").Value = '17:15'
$ws.Range("C29").Value = '21:00'
$ws.Range("D29").Value = 'BONILLA SÁNCHEZ, RAÚL FERNANDO'
$ws.Range("E29").Value = '21:15'
$ws.Range("F29").Value = '17:30'

$ws.Range("A30").Value = 'BONILLA SÁNCHEZ, RAÚL FERNANDO'
$ws.Range("B30").Value = '17:30'
$ws.Range("C30").Value = '21:15'
$ws.Range("D30").Value = 'ALVITE CORNEJO, ANGIE LUCERO'
$ws.Range("E30").Value = '21:15'
$ws.Range("F30").Value = '17:30'

$ws.Range("A31").Value = 'ALVITE CORNEJO, ANGIE LUCERO'
$ws.Range("B31").Value = '17:30'
$ws.Range("C31").Value = '21:15'
$ws.Range("D31").Value = 'MENDOZA CRUZ, LILIANA LILIANA'
$ws.Range("E31").Value = '21:30'
$ws.Range("F31").Value = '17:45'

$ws.Range("A32").Value = 'MENDOZA CRUZ, LILIANA LILIANA'
$ws.Range("B32").Value = '17:45'
$ws.Range("C32").Value = '21:30'
$ws.Range("D32").Value = 'SOTELO GONZALES, CAMILA SOFÍA'
$ws.Range("E32").Value = '21:30'
$ws.Range("F32").Value = '12:30'

$ws.Range("A33").Value = 'AYALA TAPIA, DARCIE SOL'
$ws.Range("B33").Value = '18:00'
$ws.Range("C33").Value = '21:45'
$ws.Range("D33").Value = 'CAPCHA YARANGO, DAVID'
$ws.Range("E33").Value = '21:45'
$ws.Range("F33").Value = '18:00'

$ws.Range("A34").Value = 'FLORES PAREDES, LOURDES'
$ws.Range("B34").Value = '18:00'
$ws.Range("C34").Value = '21:45'
$ws.Range("D34").Value = 'AYALA TAPIA, DARCIE SOL'
$ws.Range("E34").Value = '21:45'
$ws.Range("F34").Value = '18:00'

$ws.Range("A35").Value = 'CAPCHA YARANGO, DAVID'
$ws.Range("B35").Value = '18:00'
$ws.Range("C35").Value = '21:45'
$ws.Range("D35").Value = 'FLORES PAREDES, LOURDES'
$ws.Range("E35").Value = '21:45'
$ws.Range("F35").Value = '18:00'

$ws.Range("A36").Value = 'CARDENAS RICAPA, FABRIZIO ESTEBAN'
$ws.Range("B36").Value = '18:15'
$ws.Range("C36").Value = '22:00'
$ws.Range("D36").Value = 'CARDENAS RICAPA, FABRIZIO ESTEBAN'
$ws.Range("E36").Value = '22:00'
$ws.Range("F36").Value = '18:15'

$ws.Range("A37").Value = 'SALAS VILLANUEVA, JAMILA DASHA'
$ws.Range("B37").Value = '18:15'
$ws.Range("C37").Value = '22:00'
$ws.Range("D37").Value = 'VILCHEZ CUBA, JACK ANTHONY'
$ws.Range("E37").Value = '22:00'
$ws.Range("F37").Value = '18:15'

$ws.Range("A38").Value = 'VILCHEZ CUBA, JACK ANTHONY'
$ws.Range("B38").Value = '18:15'
$ws.Range("C38").Value = '22:00'
$ws.Range("D38").Value = 'SALAS VILLANUEVA, JAMILA DASHA'
$ws.Range("E38").Value = '22:00'
$ws.Range("F38").Value = '18:15'

$ws.Range("A39").Value = 'BRENIS LÁRTIGA, SEBASTIÁN'
$ws.Range("B39").Value = '18:30'
$ws.Range("C39").Value = '22:15'
$ws.Range("D39").Value = 'BRENIS LÁRTIGA, SEBASTIÁN'
$ws.Range("E39").Value = '22:15'
$ws.Range("F39").Value = '18:30'

$ws.Range("A40").Value = 'CORDOVA MONTALVO, MELANY KARINA'
$ws.Range("B40").Value = '19:00'
$ws.Range("C40").Value = '22:45'
$ws.Range("D40").Value = 'CORDOVA MONTALVO, MELANY KARINA'
$ws.Range("E40").Value = '22:45'
$ws.Range("F40").Value = '19:00'

$ws.Range("A41").Value = 'INGA DELGADO, CARLOS DANIEL'
$ws.Range("B41").Value = '19:00'
$ws.Range("C41").Value = '22:45'
$ws.Range("D41").Value = 'INGA DELGADO, CARLOS DANIEL'
$ws.Range("E41").Value = '22:45'
$ws.Range("F41").Value = '19:00'

$ws.Rows(42).Delete()
